$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.883.25"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").Value = "1.769.45"
$ws.Range("E3").Value = "  +2.44%  "

$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4494"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07430"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.101"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.029"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.248"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.21%  "

$ws.Range("D16").Value = "1.771.93"
$ws.Range("E16").Value = "  +2.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001062"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06431"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("E21").Value = "  +3.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.788"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.18%  "

$ws.Range("D23").Value = "27.929.75"
$ws.Range("E23").Value = "  +1.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.106"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("E27").Value = "  +3.18%  "

$ws.Range("D28").Value = "1.974.41"
$ws.Range("E28").Value = "  +2.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.162"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("E31").Value = "  +7.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09189"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.628"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.660"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("E35").Value = "  +2.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02292"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06105"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6329"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.67%  "

$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.394"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.899"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.745"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5892"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.958"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "

$ws.Range("E49").Value = "  +1.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.138"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.37%  "
